$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B3" = 2.794415499126957
    "B4" = 1.929662451756009
    "C4" = 1.79986265093306
    "D4" = -0.8854380943849716
    "C5" = -0.3103476474035083
    "D5" = -2.218321983005667
    "E5" = 1.532898100704427
    "F5" = -0.03788152406275502
    "E6" = 1.0892423430376
    "F6" = 0.3618709043640589
    "G6" = -1.675983833549222
    "H6" = 2.15807511757542
    "G7" = 1.277667191469334
    "H7" = 2.627768965131905
    "I7" = 2.659314723144024
    "J7" = 2.814292328656265
    "I8" = 1.423575078814565
    "J8" = 1.277551254953391
    "K8" = 2.906157307553836
    "L8" = 2.942581135514977
    "K9" = 4.038794034641202
    "L9" = 4.110668188518263
    "M9" = 3.586987532670949
    "N9" = 3.878230798954285
    "M10" = 5.578558913710663
    "N10" = 5.555223160690259
    "O10" = 3.42596297413984
    "P10" = 3.075158037444581
    "O11" = 2.76275821580223
    "P11" = 2.367041597905817
    "Q11" = 1.999626938280241
    "Q12" = 1.512800301290995
    "R12" = 2.590730081186199
    "S12" = 2.793289702145763
    "R13" = 2.494417544901628
    "S13" = 2.207308935472674
    "T13" = 2.821516951149361
    "U13" = 3.016819787229474
    "T14" = 3.121124374623663
    "U14" = 3.300883241600383
    "V14" = 2.535635243126988
    "W14" = 2.637488927515808
    "X14" = 2.841516658941856
    "V15" = 2.905506582474837
    "W15" = 3.099927982210238
    "X15" = 3.177336867742331
    "Y15" = 2.46048248889319
    "Z15" = 2.220932789361152
    "AA15" = 2.167670286234991
    "AB15" = 2.207173254521999
    "Y16" = 2.574315362377289
    "Z16" = 2.562407432124303
    "AA16" = 2.556000311085604
    "AB16" = 2.525788519949024
    "AC16" = 1.496881353009161
    "AD16" = 1.694296813984009
    "AE16" = 1.629686186121027
    "AF16" = 1.619782579158202
    "AC17" = 1.83915572102098
    "AD17" = 1.811820461872138
    "AE17" = 1.633300070291677
    "AF17" = 0.3071129274195616
    "AG17" = 2.257237375640031
    "AH17" = 2.059057776028594
    "AI17" = 1.937440955395164
    "AJ17" = 1.388548717051186
    "AG18" = 2.213620378726788
    "AH18" = 2.610416778758373
    "AI18" = 2.179993174715689
    "AJ18" = 0.890977499942136
    "AK18" = 2.058050235820175
    "AL18" = 2.29686889447267
    "AM18" = 3.014134262744617
    "AN18" = 1.437806261771213
    "AK19" = 2.849992723907335
    "AL19" = 3.479464952554112
    "AM19" = 4.125217580302332
    "AN19" = 2.148304186541194
    "AO19" = 2.599913004672616
    "AP19" = 2.337862417976333
    "AQ19" = 2.577692526489739
    "AR19" = -1.390622874876313
    "AO20" = 1.343460690969822
    "AP20" = 0.8329290289207147
    "AQ20" = 1.539533176834884
    "AR20" = -2.093034802586002
    "AS20" = 0.5029237023806754
    "AT20" = 0.4565833992175916
    "AU20" = 1.187924830910969
    "AV20" = 0.8461784325530575
    "AS21" = -0.7352716516441982
    "AT21" = -0.7869205535448565
    "AU21" = -0.3489999547360179
    "AV21" = -0.9869022883377543
    "AW21" = 1.005480064500386
    "AX21" = 0.5434772144153888
    "AY21" = 0.5503752294844233
    "AZ21" = 1.903706680019468
    "AW22" = 0.5134406156019233
    "AX22" = 0.03073175363270675
    "AY22" = -0.04341979710753563
    "AZ22" = 0.6808804886353492
    "BA22" = 2.403526819519342
    "BA23" = 1.910815645575914
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}